$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(31, 1).Value = '20-10-2025 10:16:01'
$ws.Cells.Item(31, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(31, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(32, 1).Value = '20-10-2025 10:19:55'
$ws.Cells.Item(32, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(32, 3).Value = 'Đăng xuất'

$ws.Cells.Item(33, 1).Value = '20-10-2025 10:27:53'
$ws.Cells.Item(33, 2).Value = 'Phạm Thị Thúy'
$ws.Cells.Item(33, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(34, 1).Value = '20-10-2025 10:28:08'
$ws.Cells.Item(34, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(34, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(35, 1).Value = '20-10-2025 10:29:01'
$ws.Cells.Item(35, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(35, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(36, 1).Value = '20-10-2025 10:30:26'
$ws.Cells.Item(36, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(36, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(37, 1).Value = '20-10-2025 10:31:24'
$ws.Cells.Item(37, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(37, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(38, 1).Value = '20-10-2025 10:32:43'
$ws.Cells.Item(38, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(38, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(39, 1).Value = '20-10-2025 10:33:46'
$ws.Cells.Item(39, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(39, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(40, 1).Value = '20-10-2025 10:34:33'
$ws.Cells.Item(40, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(40, 3).Value = 'Đăng nhập.'

$ws.Cells.Item(41, 1).Value = '20-10-2025 10:35:00'
$ws.Cells.Item(41, 2).Value = 'Nguyễn Văn Nam'
$ws.Cells.Item(41, 3).Value = 'Đăng nhập.'
